$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set the date value for A12 (2012-08-05, Excel serial 41126). Keep existing
# date number format already applied to the cell.
$ws.Range("A12").Value = 41126

# Set the text for B12 (added as a new shared string entry).
$ws.Range("B12").Value = "Implemented LibCL RadixSort, Added some keywords to lexer file, Refactored RUN macros in main"

# Update the active selection to B12 to match the saved view state.
$ws.Activate()
$ws.Range("B12").Select()
